$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the culture/group labels ("culture 1", "culture 2", "control").
# Replace the spaces in the culture labels with underscores
# ("culture 1" -> "culture_1", "culture 2" -> "culture_2").
$colA = $ws.Range("A1:A38")
$colA.Replace("culture 1", "culture_1", -4163, 1, $false, $false, $true)
$colA.Replace("culture 2", "culture_2", -4163, 1, $false, $false, $true)
